$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6.403577508329116
$ws.Range("G2").Value = 88.56999999999999
$ws.Range("H2").Value = 0.03183715456471803
$ws.Range("J2").Value = 0.01843715456471803
$ws.Range("K2").Value = 0.03663715456471803
$ws.Range("F3").Value = 5.411796627969051
$ws.Range("G3").Value = 95.72
$ws.Range("H3").Value = 0.03553937432865162
$ws.Range("J3").Value = 0.02323937432865162
$ws.Range("K3").Value = 0.04143937432865162
$ws.Range("F4").Value = 11.42001574762284
$ws.Range("G4").Value = 76.84
$ws.Range("H4").Value = 0.03147061866053491
$ws.Range("I4").Value = 0.0191
$ws.Range("J4").Value = 0.01237061866053491
$ws.Range("K4").Value = 0.03057061866053491
$ws.Range("F5").Value = 6.296728015749271
$ws.Range("G5").Value = 83.78
$ws.Range("H5").Value = 0.02850470550880613
$ws.Range("J5").Value = 0.01510470550880613
$ws.Range("K5").Value = 0.03330470550880613
$ws.Range("F6").Value = 9.483029329420187
$ws.Range("G6").Value = 78.89
$ws.Range("H6").Value = 0.03165738164471943
$ws.Range("I6").Value = 0.0169
$ws.Range("J6").Value = 0.01475738164471943
$ws.Range("K6").Value = 0.03295738164471943
$ws.Range("F7").Value = 7.803029445818093
$ws.Range("G7").Value = 84.75999999999999
$ws.Range("H7").Value = 0.03180184694148602
$ws.Range("I7").Value = 0.01502
$ws.Range("J7").Value = 0.01678184694148602
$ws.Range("K7").Value = 0.03498184694148603
